$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.736.51'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.26%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.978.10'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +4.71%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '484.17'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +9.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.27'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.07%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.735'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.170'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +12.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000361'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +16.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.42'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.605.10'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.50'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.06'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.31%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.980.59'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +5.75%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.09'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.16'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '67.750.71'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '436.73'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +5.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.44'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +5.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.59'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.13'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.66'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +7.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '38.99'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +5.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.96'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.11'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '725.63'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.93%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.131'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.93%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.40'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.28%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.90%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '42.24'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0882'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +30.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '59.68'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.59%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.15%  '
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.38'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.35%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.05'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +6.24%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.26'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +7.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.342'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +8.31%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.47'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.53'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.71%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '148.93'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.78%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.22'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.09%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.01%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.06%  '
